{"js": "// Remove the trailing \"Ver no Jupiter / Save pdf / docx\" helper paragraph,\n// the \"\u00a9 2020 ... Jekyll ...\" footer paragraph, and the blank paragraph that\n// separates them from the preceding \"LOB1019: F\u00edsica II (Requisito fraco)\"\n// paragraph \u2014 i.e. collapse the four paragraphs down to just the first one.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetSnippets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\n\nconst items = paragraphs.items;\n\n// Find the index of the paragraph that marks the requirement line; the\n// blank paragraph right after it (before the two footer paragraphs) should\n// be deleted along with the footer paragraphs themselves.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOB1019: F\u00edsica II (Requisito fraco)\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (anchorIndex !== -1 && anchorIndex + 1 < items.length && items[anchorIndex + 1].text === \"\") {\n  toDelete.push(items[anchorIndex + 1]);\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (targetSnippets.some((snippet) => text.indexOf(snippet) !== -1)) {\n    toDelete.push(items[i]);\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter / Save pdf / docx\" helper paragraph,\n# the \"(c) 2020 ... Jekyll ...\" footer paragraph, and the blank paragraph\n# that separates them from the preceding\n# \"LOB1019: Fisica II (Requisito fraco)\" paragraph.\n\n$d = $word.ActiveDocument\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -like \"*LOB1019*Requisito fraco*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\n$indicesToDelete = @()\n\nif ($anchorIndex -ne -1) {\n    $blankIndex = $anchorIndex + 1\n    if ($blankIndex -le $d.Paragraphs.Count) {\n        $blankText = $d.Paragraphs($blankIndex).Range.Text.Trim()\n        if ($blankText -eq \"\") {\n            $indicesToDelete += $blankIndex\n        }\n    }\n}\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if (($t -like \"*Ver no Jupiter*\") -or ($t -like \"*Jekyll*\")) {\n        $indicesToDelete += $i\n    }\n}\n\n# Delete from the highest index down so earlier indices stay valid.\n$indicesToDelete = $indicesToDelete | Sort-Object -Descending -Unique\nforeach ($idx in $indicesToDelete) {\n    $d.Paragraphs($idx).Range.Delete()\n}\n"}
